$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 128, pushing all existing rows (128-212)
# down by one (to 129-213), matching the target dimension A1:R213.
$ws.Rows.Item(128).Insert()

$ws.Range("A128").Value = 6
$ws.Range("B128").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C128").Value = 'Metropolitana'
$ws.Range("D128").Value = 44673
$ws.Range("E128").Value = 13
$ws.Range("F128").Value = 100112022
$ws.Range("G128").Value = 'Arveja Verde'
$ws.Range("H128").Value = 'Perfection'
$ws.Range("I128").Value = 'Primera'
$ws.Range("J128").Value = 90
$ws.Range("K128").Value = 27000
$ws.Range("L128").Value = 30000
$ws.Range("M128").Value = 28333
$ws.Range("N128").Value = '$/malla 25 kilos'
$ws.Range("O128").Value = 'Provincia de Huasco'
$ws.Range("P128").Value = 1133
$ws.Range("Q128").Value = 25
$ws.Range("R128").Value = 'Hortaliza'
